$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old extra rows (3-5) so the sheet keeps only the header and one
# autocomplete-ready data row, matching the new A1:E2 used range.
$ws.Range("A3:E5").Clear()

# C2 and E2 must stay textual ("1" and "b"), not be auto-coerced to numbers,
# so mark them as Text before writing the values.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

# Update row 2 with the new sample values.
$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "b@b"
$ws.Range("E2").Value = "b"
